$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (characters, customWidth) ---
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334
$ws.Columns.Item(3).ColumnWidth = 16.0
$ws.Columns.Item(4).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 16.0
$ws.Columns.Item(6).ColumnWidth = 10.0
$ws.Columns.Item(7).ColumnWidth = 10.0
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
$ws.Columns.Item(9).ColumnWidth = 10.0
$ws.Columns.Item(10).ColumnWidth = 13.5

# --- Header row (now wrapped onto multiple lines + updated weights) ---
$ws.Range("B1").Value = "Стоимость`nТС`n(0.328)"
$ws.Range("C1").Value = "Стоимость`nобслуживания`nТС`n(0.159)"
$ws.Range("D1").Value = "Надёжность`n(0.232)"
$ws.Range("E1").Value = "Безопасность`n(0.107)"
$ws.Range("F1").Value = "Дизайн`n(0.048)"
$ws.Range("G1").Value = "Комфорт`n(0.071)"
$ws.Range("H1").Value = "Мощность`n(0.033)"
$ws.Range("I1").Value = "Год`nвыпуска`n(0.023)"
$ws.Range("J1").Value = "Глобальные`nприоритеты`nвыбора"

# --- Row 2: Kia Rio ---
$ws.Range("A2").Value = "Kia`nRio"
$ws.Range("B2").Value = "'0.276"
$ws.Range("C2").Value = "'0.179"
$ws.Range("D2").Value = "'0.097"
$ws.Range("I2").Value = "'0.222"
$ws.Range("J2").Value = "'0.166"

# --- Row 3: Volkswagen Golf ---
$ws.Range("A3").Value = "Volkswagen`nGolf"
$ws.Range("B3").Value = "'0.092"
$ws.Range("C3").Value = "'0.101"
$ws.Range("F3").Value = "'0.228"
$ws.Range("H3").Value = "'0.228"
$ws.Range("I3").Value = "'0.092"
$ws.Range("J3").Value = "'0.104"

# --- Row 4: Toyota Corolla ---
$ws.Range("A4").Value = "Toyota`nCorolla"
$ws.Range("B4").Value = "'0.157"
$ws.Range("C4").Value = "'0.316"
$ws.Range("D4").Value = "'0.249"
$ws.Range("E4").Value = "'0.228"
$ws.Range("G4").Value = "'0.228"
$ws.Range("J4").Value = "'0.218"

# --- Row 5: Skoda Octavia ---
$ws.Range("A5").Value = "Skoda`nOctavia"
$ws.Range("B5").Value = "'0.157"
$ws.Range("C5").Value = "'0.179"
$ws.Range("D5").Value = "'0.157"
$ws.Range("I5").Value = "'0.222"
$ws.Range("J5").Value = "'0.155"

# --- Row 6: BMW 3 Series ---
$ws.Range("A6").Value = "BMW`n3`nSeries"
$ws.Range("B6").Value = "'0.043"
$ws.Range("C6").Value = "'0.045"
$ws.Range("E6").Value = "'0.362"
$ws.Range("F6").Value = "'0.362"
$ws.Range("G6").Value = "'0.362"
$ws.Range("H6").Value = "'0.362"
$ws.Range("J6").Value = "'0.204"

# --- Row 7: Hyundai Solaris ---
$ws.Range("A7").Value = "Hyundai`nSolaris"
$ws.Range("B7").Value = "'0.276"
$ws.Range("C7").Value = "'0.179"
$ws.Range("I7").Value = "'0.092"
$ws.Range("J7").Value = "'0.154"

# The cells above now wrap onto multiple lines; restore each row to its
# natural (default) height instead of leaving an explicit custom height.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()
